$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.043.32"
$ws.Range("E2").Value = "  -1.72%  "
$ws.Range("D3").Value = "2.566.87"
$ws.Range("E3").Value = "  -3.39%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.87"
$ws.Range("E5").Value = "  -1.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.70"
$ws.Range("E6").Value = "  -4.52%  "
$ws.Range("E7").Value = "  -2.55%  "
$ws.Range("E8").Value = "  -0.22%  "
$ws.Range("E9").Value = "  -3.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.32"
$ws.Range("E10").Value = "  -3.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0813"
$ws.Range("E11").Value = "  -1.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.69"
$ws.Range("E12").Value = "  -2.82%  "
$ws.Range("E13").Value = "  +6.90%  "
$ws.Range("D14").Value = "2.961.77"
$ws.Range("E14").Value = "  -4.03%  "
$ws.Range("D15").Value = "2.532.39"
$ws.Range("E15").Value = "  -5.42%  "
$ws.Range("E16").Value = "  -2.69%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.24"
$ws.Range("E17").Value = "  -3.77%  "
$ws.Range("D18").Value = "43.112.88"
$ws.Range("E18").Value = "  -2.20%  "
$ws.Range("D19").Value = "0.0₃0987"
$ws.Range("E19").Value = "  -0.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.69"
$ws.Range("E20").Value = "  +0.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.58"
$ws.Range("E21").Value = "  -3.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.82"
$ws.Range("E22").Value = "  -3.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "254.28"
$ws.Range("E23").Value = "  -6.90%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.95"
$ws.Range("E24").Value = "  -1.41%  "
$ws.Range("E25").Value = "  -6.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "29.05"
$ws.Range("E26").Value = "  -3.90%  "
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("E28").Value = "  -1.73%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "37.27"
$ws.Range("E29").Value = "  -0.86%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.14"
$ws.Range("E30").Value = "  -4.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.11"
$ws.Range("E31").Value = "  +0.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "153.75"
$ws.Range("E32").Value = "  +0.22%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.77"
$ws.Range("E33").Value = "  -1.33%  "
$ws.Range("E34").Value = "  -7.80%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.14"
$ws.Range("E35").Value = "  -7.61%  "
$ws.Range("E36").Value = "  -3.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.113"
$ws.Range("E37").Value = "  -4.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "17.71"
$ws.Range("E38").Value = "  +11.94%  "
$ws.Range("E39").Value = "  -2.19%  "
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "22.87"
$ws.Range("E40").Value = "  -8.84%  "
$ws.Range("B41").Value = "ApeXProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.19"
$ws.Range("E41").Value = "  +36.30%  "
$ws.Range("B42").Value = "NEARProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.42"
$ws.Range("E42").Value = "  -3.83%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0311"
$ws.Range("E43").Value = "  -2.89%  "
$ws.Range("E44").Value = "  -0.60%  "
$ws.Range("D45").Value = "2.101.14"
$ws.Range("E45").Value = "  -0.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.998"
$ws.Range("E46").Value = "  -0.18%  "
$ws.Range("E47").Value = "  +0.84%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.76"
$ws.Range("E48").Value = "  -5.44%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "106.09"
$ws.Range("E49").Value = "  -2.37%  "
$ws.Range("D50").Value = "2.817.20"
$ws.Range("E50").Value = "  -3.89%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "74.89"
$ws.Range("E51").Value = "  +5.84%  "
